$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.532.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.188.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.398"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.715"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.57%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.186.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.568"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.61%  "

$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91.025.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.766.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.179.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000213"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +60.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "438.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.344.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.156"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +32.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "519.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "171.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.743"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.610"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.81%  "
